$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.542.91"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.058.59"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "386.47"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.09"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.78"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.547.01"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.59"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.78"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.059.81"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.973"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.68"
$ws.Range("E18").Value = "  -4.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.595.32"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.45"
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.18"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.86"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.15"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("E26").Value = "  +4.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.84"
$ws.Range("E27").Value = "  +2.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.30"
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.28"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.78"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0450"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "49.94"
$ws.Range("E36").Value = "  -3.19%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.293"
$ws.Range("E39").Value = "  +8.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.92"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.17"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.74"
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.10"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.033.79"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.360.91"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.207"
$ws.Range("E51").Value = "  +7.06%  "
